$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet ("ODI Bowling")
# so it lands at the end of the tab order as "ODI Batting Extra".
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "ODI Batting Extra"

# Header row (copy the header formatting used by the other sheets
# so the new header row matches the existing bold/bordered/centered style).
$headerSrc = $wb.Worksheets.Item("Player Info")
$headerSrc.Range("A1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "BATTING_POSITION"
$ws.Range("C1").Value = "NUM_4"
$ws.Range("D1").Value = "NUM_6"
$ws.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Range("F1").Value = "MAN_OF_MATCH"

# Row 2
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "4711"
$ws.Range("B2").Value = 6
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "3"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "12.26%"
$ws.Range("F2").Value = "NO"

# Row 3
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "4713"
$ws.Range("B3").Value = 6
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "0"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "0"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.31%"
$ws.Range("F3").Value = "NO"
